$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-20 18:52:58"
$wsZh.Range("H2").Value = "2016-03-20 18:53:18"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-20 18:53:01"
$wsDe.Range("H2").Value = "2016-03-20 18:53:23"
